$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.035739180850205
$ws.Range("D2").Value = 1.0365018361634
$ws.Range("E2").Value = 1.04407915353227
$ws.Range("F2").Value = 1.053369209529297
$ws.Range("I2").Value = 1.037041961587844
$ws.Range("J2").Value = 1.040851311308157
$ws.Range("K2").Value = 1.039295395551752
$ws.Range("L2").Value = 1.046851213677209
$ws.Range("M2").Value = 1.056115368416734
$ws.Range("N2").Value = 1.042329439557529

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.036705795916226
$ws.Range("D3").Value = 1.037335210576458
$ws.Range("E3").Value = 1.044945200196183
$ws.Range("F3").Value = 1.054338083335344
$ws.Range("I3").Value = 1.037212169769669
$ws.Range("J3").Value = 1.041461513686232
$ws.Range("K3").Value = 1.039938321567864
$ws.Range("L3").Value = 1.047528282162623
$ws.Range("M3").Value = 1.056896873952488
$ws.Range("N3").Value = 1.042940508492972

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.0373313370646
$ws.Range("D4").Value = 1.037874822328026
$ws.Range("E4").Value = 1.045506022078144
$ws.Range("F4").Value = 1.054965428013442
$ws.Range("I4").Value = 1.037320362382317
$ws.Range("J4").Value = 1.041855826037568
$ws.Range("K4").Value = 1.040354055020022
$ws.Range("L4").Value = 1.047966182353089
$ws.Range("M4").Value = 1.057402358624128
$ws.Range("N4").Value = 1.043335380813075

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.037594332136853
$ws.Range("D5").Value = 1.038101760785552
$ws.Range("E5").Value = 1.045741893795702
$ws.Range("F5").Value = 1.055229262738761
$ws.Range("I5").Value = 1.037365380624386
$ws.Range("J5").Value = 1.042021467335427
$ws.Range("K5").Value = 1.040528760762698
$ws.Range("L5").Value = 1.048150224752925
$ws.Range("M5").Value = 1.057614815036017
$ws.Range("N5").Value = 1.043501257340576

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.037638491175217
$ws.Range("D6").Value = 1.038139869741096
$ws.Range("E6").Value = 1.045781503653808
$ws.Range("F6").Value = 1.05527356752252
$ws.Range("I6").Value = 1.03737291204343
$ws.Range("J6").Value = 1.042049271729208
$ws.Range("K6").Value = 1.040558090575781
$ws.Range("L6").Value = 1.048181123272672
$ws.Range("M6").Value = 1.057650484477825
$ws.Range("N6").Value = 1.043529101219786

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.037334851149526
$ws.Range("D7").Value = 1.037877854355632
$ws.Range("E7").Value = 1.045509173407487
$ws.Range("F7").Value = 1.054968952998349
$ws.Range("I7").Value = 1.037320965749909
$ws.Range("J7").Value = 1.041858039847246
$ws.Range("K7").Value = 1.040356389717435
$ws.Range("L7").Value = 1.047968641737613
$ws.Range("M7").Value = 1.057405197672338
$ws.Range("N7").Value = 1.043337597766617

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.036065836630401
$ws.Range("D8").Value = 1.036783403128507
$ws.Range("E8").Value = 1.044371747952245
$ws.Range("F8").Value = 1.053696557601435
$ws.Range("I8").Value = 1.037099886249488
$ws.Range("J8").Value = 1.041057641146594
$ws.Range("K8").Value = 1.03951273338143
$ws.Range("L8").Value = 1.047080074721404
$ws.Range("M8").Value = 1.056379522989992
$ws.Range("N8").Value = 1.042536062408003

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.033830290233066
$ws.Range("D9").Value = 1.03485767009201
$ws.Range("E9").Value = 1.042370817174139
$ws.Range("F9").Value = 1.051457691338443
$ws.Range("I9").Value = 1.036695460239265
$ws.Range("J9").Value = 1.039643220207776
$ws.Range("K9").Value = 1.03802397880524
$ws.Range("L9").Value = 1.045512751630305
$ws.Range("M9").Value = 1.05457065516587
$ws.Range("N9").Value = 1.041119632829215

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.032340384432139
$ws.Range("D10").Value = 1.033575818487133
$ws.Range("E10").Value = 1.041039188759206
$ws.Range("F10").Value = 1.049967371180601
$ws.Range("I10").Value = 1.036415892769166
$ws.Range("J10").Value = 1.038697622859061
$ws.Range("K10").Value = 1.037030100794552
$ws.Range("L10").Value = 1.044466878996869
$ws.Range("M10").Value = 1.053363797167409
$ws.Range("N10").Value = 1.040172692623804

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.031695355636925
$ws.Range("D11").Value = 1.033021243768894
$ws.Range("E11").Value = 1.040463144276652
$ws.Range("F11").Value = 1.049322595159394
$ws.Range("I11").Value = 1.036292484614015
$ws.Range("J11").Value = 1.038287550315351
$ws.Range("K11").Value = 1.036599425310239
$ws.Range("L11").Value = 1.044013780690634
$ws.Range("M11").Value = 1.052841002255492
$ws.Range("N11").Value = 1.039762037730045

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.031455780352453
$ws.Range("D12").Value = 1.032815322526435
$ws.Range("E12").Value = 1.04024926080955
$ws.Range("F12").Value = 1.049083179167663
$ws.Range("I12").Value = 1.036246292293345
$ws.Range("J12").Value = 1.038135138322587
$ws.Range("K12").Value = 1.036439406036276
$ws.Range("L12").Value = 1.043845446302273
$ws.Range("M12").Value = 1.052646781210493
$ws.Range("N12").Value = 1.039609409294769

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.031507169279988
$ws.Range("D13").Value = 1.032859490034137
$ws.Range("E13").Value = 1.040295135671314
$ws.Range("F13").Value = 1.04913453095932
$ws.Range("I13").Value = 1.036256216669601
$ws.Range("J13").Value = 1.03816783539935
$ws.Range("K13").Value = 1.036473732837651
$ws.Range("L13").Value = 1.043881556102741
$ws.Range("M13").Value = 1.052688443705245
$ws.Range("N13").Value = 1.039642152805132

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.031675551916459
$ws.Range("D14").Value = 1.033004220769636
$ws.Range("E14").Value = 1.040445462860386
$ws.Range("F14").Value = 1.049302803269249
$ws.Range("I14").Value = 1.036288673540985
$ws.Range("J14").Value = 1.038274953781055
$ws.Range("K14").Value = 1.036586199025778
$ws.Range("L14").Value = 1.043999866796226
$ws.Range("M14").Value = 1.052824948525396
$ws.Range("N14").Value = 1.039749423307226

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.031779300302919
$ws.Range("D15").Value = 1.033093403791228
$ws.Range("E15").Value = 1.040538095711393
$ws.Range("F15").Value = 1.049406492360161
$ws.Range("I15").Value = 1.036308624526601
$ws.Range("J15").Value = 1.038340940678097
$ws.Range("K15").Value = 1.036655486922899
$ws.Range("L15").Value = 1.044072757504612
$ws.Range("M15").Value = 1.052909049472551
$ws.Range("N15").Value = 1.03981550391323

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.032383195240121
$ws.Range("D16").Value = 1.033612633893609
$ws.Range("E16").Value = 1.04107743081081
$ws.Range("F16").Value = 1.050010174358726
$ws.Range("I16").Value = 1.036424033422691
$ws.Range("J16").Value = 1.038724824969222
$ws.Range("K16").Value = 1.037058676641027
$ws.Range("L16").Value = 1.044496944895164
$ws.Range("M16").Value = 1.053398488842313
$ws.Range("N16").Value = 1.040199933364081

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.032762032496762
$ws.Range("D17").Value = 1.03393846127282
$ws.Range("E17").Value = 1.041415891987213
$ws.Range("F17").Value = 1.050388994086999
$ws.Range("I17").Value = 1.036495796649575
$ws.Range("J17").Value = 1.038965459275808
$ws.Range("K17").Value = 1.037311501870662
$ws.Range("L17").Value = 1.04476296583849
$ws.Range("M17").Value = 1.053705443672811
$ws.Range("N17").Value = 1.040440909398983

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.032983012458669
$ws.Range("D18").Value = 1.034128556676972
$ws.Range("E18").Value = 1.041613364569297
$ws.Range("F18").Value = 1.050610005493983
$ws.Range("I18").Value = 1.036537427800284
$ws.Range("J18").Value = 1.039105757068411
$ws.Range("K18").Value = 1.037458939644289
$ws.Range("L18").Value = 1.044918109246142
$ws.Range("M18").Value = 1.05388446419083
$ws.Range("N18").Value = 1.040581406430545

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.033058362695643
$ws.Range("D19").Value = 1.034193382034394
$ws.Range("E19").Value = 1.041680706727482
$ws.Range("F19").Value = 1.050685373506405
$ws.Range("I19").Value = 1.036551584405185
$ws.Range("J19").Value = 1.039153584768415
$ws.Range("K19").Value = 1.037509206848009
$ws.Range("L19").Value = 1.044971005363552
$ws.Range("M19").Value = 1.053945501963828
$ws.Range("N19").Value = 1.040629302051369

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.032721385743862
$ws.Range("D20").Value = 1.033903498332911
$ws.Range("E20").Value = 1.041379572722886
$ws.Range("F20").Value = 1.050348344901326
$ws.Range("I20").Value = 1.036488120623018
$ws.Range("J20").Value = 1.038939647705777
$ws.Range("K20").Value = 1.037284379306975
$ws.Range("L20").Value = 1.044734426589501
$ws.Range("M20").Value = 1.053672512494107
$ws.Range("N20").Value = 1.040415061173562

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.031625966949433
$ws.Range("D21").Value = 1.032961599165466
$ws.Range("E21").Value = 1.040401192892309
$ws.Range("F21").Value = 1.049253248981493
$ws.Range("I21").Value = 1.036279125542375
$ws.Range("J21").Value = 1.038243412642958
$ws.Range("K21").Value = 1.036553081833589
$ws.Range("L21").Value = 1.043965028163984
$ws.Range("M21").Value = 1.052784752126815
$ws.Range("N21").Value = 1.039717837377095

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.030937333064365
$ws.Range("D22").Value = 1.032369809762647
$ws.Range("E22").Value = 1.039786539555512
$ws.Range("F22").Value = 1.048565196721858
$ws.Range("I22").Value = 1.036145679752085
$ws.Range("J22").Value = 1.037805125532707
$ws.Range("K22").Value = 1.036093013268071
$ws.Range("L22").Value = 1.043481083407222
$ws.Range("M22").Value = 1.05222639878361
$ws.Range("N22").Value = 1.039278927848875

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.031302381306732
$ws.Range("D23").Value = 1.032683488358214
$ws.Range("E23").Value = 1.040112331846509
$ws.Range("F23").Value = 1.048929900588774
$ws.Range("I23").Value = 1.036216615288676
$ws.Range("J23").Value = 1.038037520424965
$ws.Range("K23").Value = 1.036336929908521
$ws.Range("L23").Value = 1.043737649736694
$ws.Range("M23").Value = 1.052522409507989
$ws.Range("N23").Value = 1.039511652768536

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.032739752229089
$ws.Range("D24").Value = 1.033919296439025
$ws.Range("E24").Value = 1.041395983668598
$ws.Range("F24").Value = 1.050366712357149
$ws.Range("I24").Value = 1.036491589790811
$ws.Range("J24").Value = 1.038951311029034
$ws.Range("K24").Value = 1.037296634921555
$ws.Range("L24").Value = 1.044747322316006
$ws.Range("M24").Value = 1.053687392741306
$ws.Range("N24").Value = 1.040426741060076

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.034408154824133
$ws.Range("D25").Value = 1.03535517588259
$ws.Range("E25").Value = 1.04288770079502
$ws.Range("F25").Value = 1.05203609898308
$ws.Range("I25").Value = 1.036801770595813
$ws.Range("J25").Value = 1.040009352669185
$ws.Range("K25").Value = 1.038409103778471
$ws.Range("L25").Value = 1.045918120789566
$ws.Range("M25").Value = 1.05503846206426
$ws.Range("N25").Value = 1.041486285240715
